# Revert "Added a style":
#  - remove the custom "MyStyle" paragraph style reference from the first
#    paragraph (drop the <w:pPr><w:pStyle .../></w:pPr> entirely)
#  - move the _GoBack bookmark (bookmarkStart/bookmarkEnd) so it comes
#    after the run instead of before it
#  - delete the now-unused "My Style" style definition from styles.xml

$d = $word.ActiveDocument

$para = $d.Paragraphs(1)
$range = $para.Range

# Rebuild the paragraph's raw XML without the pPr/pStyle override and with
# the bookmark markers moved to after the text run. Using InsertXML (rather
# than Range.Text / Style assignment) lets us drop the <w:pPr> element
# completely instead of leaving an explicit <w:pStyle w:val="Normal"/>.
$newParaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>Hello world</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"
$range.InsertXML($newParaXml)

# Drop the custom style definition added by the reverted commit.
$d.Styles("My Style").Delete()
